# Scheduled-runner price refresh: update currentAveragePrice / LevePrice / LeveProfit
# columns (H, I, J, K, L, M, N) for the affected Leve rows on each job sheet.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 11
$ws.Range("H11").Value = 43
$ws.Range("I11").Value = 43
$ws.Range("K11").Value = 43
$ws.Range("M11").Value = 97
# Row 17
$ws.Range("H17").Value = 1814.7
$ws.Range("J17").Value = 1898.5227
$ws.Range("L17").Value = 5695.5681
$ws.Range("N17").Value = -6031.5681
# Row 125
$ws.Range("H125").Value = 3996.25
$ws.Range("I125").Value = 0
$ws.Range("K125").Value = 0
$ws.Range("M125").ClearContents()
# Row 131
$ws.Range("H131").Value = 3456.8667
$ws.Range("I131").Value = 1200.1111
$ws.Range("J131").Value = 6842
$ws.Range("K131").Value = 3600.3333
$ws.Range("L131").Value = 20526
$ws.Range("M131").Value = 1439.6667
$ws.Range("N131").Value = -30606
# Row 138
$ws.Range("H138").Value = 3081.6333
$ws.Range("I138").Value = 2309.5
$ws.Range("J138").Value = 4239.8335
$ws.Range("K138").Value = 6928.5
$ws.Range("L138").Value = 12719.5005
$ws.Range("M138").Value = -1788.5
$ws.Range("N138").Value = -22999.5005
# Row 141
$ws.Range("H141").Value = 1613.5
$ws.Range("J141").Value = 2105
$ws.Range("L141").Value = 6315
$ws.Range("N141").Value = -16675

$ws = $wb.Worksheets.Item("ARM")
# Row 5
$ws.Range("H5").Value = 324.875
$ws.Range("J5").Value = 482.33334
$ws.Range("L5").Value = 482.33334
$ws.Range("N5").Value = -706.33334
# Row 32
$ws.Range("H32").Value = 3371562.5
$ws.Range("I32").Value = 4331912.5
$ws.Range("K32").Value = 4331912.5
$ws.Range("M32").Value = -4331625.5
# Row 61
$ws.Range("H61").Value = 1598462.8
$ws.Range("I61").Value = 1677535.9
$ws.Range("K61").Value = 1677535.9
$ws.Range("M61").Value = -1677323.9
# Row 110
$ws.Range("H110").Value = 1254.1333
$ws.Range("I110").Value = 628.5
$ws.Range("J110").Value = 10013
$ws.Range("K110").Value = 628.5
$ws.Range("L110").Value = 10013
$ws.Range("M110").Value = 1416.5
$ws.Range("N110").Value = -14103
# Row 122
$ws.Range("H122").Value = 3750
$ws.Range("I122").Value = 3000
$ws.Range("K122").Value = 9000
$ws.Range("M122").Value = -6550
# Row 132
$ws.Range("H132").Value = 532046.6
$ws.Range("I132").Value = 584081.3
$ws.Range("K132").Value = 1752243.9
$ws.Range("M132").Value = -1749713.9
# Row 136
$ws.Range("H136").Value = 1598462.8
$ws.Range("I136").Value = 1677535.9
$ws.Range("K136").Value = 5032607.699999999
$ws.Range("M136").Value = -5030057.699999999

$ws = $wb.Worksheets.Item("BSM")
# Row 4
$ws.Range("H4").Value = 324.875
$ws.Range("J4").Value = 482.33334
$ws.Range("L4").Value = 482.33334
$ws.Range("N4").Value = -712.33334
# Row 17
$ws.Range("H17").Value = 1000
$ws.Range("J17").Value = 1000
$ws.Range("L17").Value = 1000
$ws.Range("N17").Value = -1344
# Row 22
$ws.Range("H22").Value = 1409
$ws.Range("I22").Value = 1577.6666
$ws.Range("J22").Value = 650
$ws.Range("K22").Value = 1577.6666
$ws.Range("L22").Value = 650
$ws.Range("M22").Value = -1404.6666
$ws.Range("N22").Value = -996
# Row 105
$ws.Range("H105").Value = 2124.0625
$ws.Range("I105").Value = 2019.1538
$ws.Range("K105").Value = 2019.1538
$ws.Range("M105").Value = -272.1538
# Row 107
$ws.Range("H107").Value = 1963.0416
$ws.Range("I107").Value = 2046
$ws.Range("K107").Value = 2046
$ws.Range("M107").Value = -126

$ws = $wb.Worksheets.Item("CRP")
# Row 7
$ws.Range("H7").Value = 228.66667
$ws.Range("I7").Value = 43.4
$ws.Range("J7").Value = 361
$ws.Range("K7").Value = 43.4
$ws.Range("L7").Value = 361
$ws.Range("M7").Value = 69.59999999999999
$ws.Range("N7").Value = -587
# Row 22
$ws.Range("H22").Value = 924.3570999999999
$ws.Range("J22").Value = 1511.75
$ws.Range("L22").Value = 1511.75
$ws.Range("N22").Value = -2211.75
# Row 31
$ws.Range("H31").Value = 5343.76
$ws.Range("I31").Value = 2080.0715
$ws.Range("J31").Value = 9497.546
$ws.Range("K31").Value = 2080.0715
$ws.Range("L31").Value = 9497.546
$ws.Range("M31").Value = -1785.0715
$ws.Range("N31").Value = -10087.546
# Row 34
$ws.Range("H34").Value = 5343.76
$ws.Range("I34").Value = 2080.0715
$ws.Range("J34").Value = 9497.546
$ws.Range("K34").Value = 2080.0715
$ws.Range("L34").Value = 9497.546
$ws.Range("M34").Value = -1878.0715
$ws.Range("N34").Value = -9901.546
# Row 43
$ws.Range("H43").Value = 85547.5
$ws.Range("I43").Value = 0
$ws.Range("K43").Value = 0
$ws.Range("M43").ClearContents()
# Row 58
$ws.Range("H58").Value = 825582.0600000001
$ws.Range("I58").Value = 1123958.2
$ws.Range("J58").Value = 5047.5
$ws.Range("K58").Value = 1123958.2
$ws.Range("L58").Value = 5047.5
$ws.Range("M58").Value = -1123755.2
$ws.Range("N58").Value = -5453.5
# Row 101
$ws.Range("H101").Value = 85547.5
$ws.Range("I101").Value = 0
$ws.Range("K101").Value = 0
$ws.Range("M101").ClearContents()
# Row 132
$ws.Range("H132").Value = 12649.667
$ws.Range("I132").Value = 13338.964
$ws.Range("J132").Value = 2999.5
$ws.Range("K132").Value = 40016.892
$ws.Range("L132").Value = 8998.5
$ws.Range("M132").Value = -37486.892
$ws.Range("N132").Value = -14058.5
# Row 134
$ws.Range("H134").Value = 1480.1428
$ws.Range("I134").Value = 1513.88
$ws.Range("K134").Value = 4541.64
$ws.Range("M134").Value = -2006.64
# Row 135
$ws.Range("H135").Value = 99000
$ws.Range("J135").Value = 99000
$ws.Range("L135").Value = 99000
$ws.Range("N135").Value = -109140
# Row 136
$ws.Range("H136").Value = 825582.0600000001
$ws.Range("I136").Value = 1123958.2
$ws.Range("J136").Value = 5047.5
$ws.Range("K136").Value = 3371874.6
$ws.Range("L136").Value = 15142.5
$ws.Range("M136").Value = -3369324.6
$ws.Range("N136").Value = -20242.5

$ws = $wb.Worksheets.Item("CUL")
# Row 3
$ws.Range("H3").Value = 5298.579
$ws.Range("I3").Value = 1257.8667
$ws.Range("J3").Value = 20451.25
$ws.Range("K3").Value = 3773.6001
$ws.Range("L3").Value = 61353.75
$ws.Range("M3").Value = -3661.6001
$ws.Range("N3").Value = -61577.75
# Row 137
$ws.Range("H137").Value = 4025.077
$ws.Range("I137").Value = 2435.1
$ws.Range("J137").Value = 9325
$ws.Range("K137").Value = 7305.299999999999
$ws.Range("L137").Value = 27975
$ws.Range("M137").Value = -2205.299999999999
$ws.Range("N137").Value = -38175

$ws = $wb.Worksheets.Item("GSM")
# Row 39
$ws.Range("H39").Value = 26704.2
$ws.Range("J39").Value = 26704.2
$ws.Range("L39").Value = 26704.2
$ws.Range("N39").Value = -27768.2
# Row 106
$ws.Range("H106").Value = 150000
$ws.Range("J106").Value = 150000
$ws.Range("L106").Value = 150000
$ws.Range("N106").Value = -152524
# Row 126
$ws.Range("H126").Value = 797961.6
$ws.Range("I126").Value = 1853939.6
$ws.Range("J126").Value = 5978.0835
$ws.Range("K126").Value = 5561818.800000001
$ws.Range("L126").Value = 17934.2505
$ws.Range("M126").Value = -5559348.800000001
$ws.Range("N126").Value = -22874.2505
# Row 132
$ws.Range("H132").Value = 2406.6365
$ws.Range("I132").Value = 1436.75
$ws.Range("J132").Value = 4993
$ws.Range("K132").Value = 4310.25
$ws.Range("L132").Value = 14979
$ws.Range("M132").Value = -1780.25
$ws.Range("N132").Value = -20039

$ws = $wb.Worksheets.Item("LTW")
# Row 58
$ws.Range("H58").Value = 0
$ws.Range("I58").Value = 0
$ws.Range("J58").Value = 0
$ws.Range("K58").Value = 0
$ws.Range("L58").Value = 0
$ws.Range("M58").ClearContents()
$ws.Range("N58").ClearContents()
# Row 101
$ws.Range("H101").Value = 65744.60000000001
$ws.Range("J101").Value = 65744.60000000001
$ws.Range("L101").Value = 65744.60000000001
$ws.Range("N101").Value = -72234.60000000001
# Row 132
$ws.Range("H132").Value = 3274.8865
$ws.Range("I132").Value = 2936.7104
$ws.Range("J132").Value = 5416.6665
$ws.Range("K132").Value = 8810.1312
$ws.Range("L132").Value = 16249.9995
$ws.Range("M132").Value = -6280.1312
$ws.Range("N132").Value = -21309.9995

$ws = $wb.Worksheets.Item("WVR")
# Row 9
$ws.Range("H9").Value = 5000
$ws.Range("I9").Value = 0
$ws.Range("J9").Value = 5000
$ws.Range("K9").Value = 0
$ws.Range("L9").Value = 5000
$ws.Range("M9").ClearContents()
$ws.Range("N9").Value = -5280
# Row 16
$ws.Range("H16").Value = 71000
$ws.Range("J16").Value = 71000
$ws.Range("L16").Value = 71000
$ws.Range("N16").Value = -71584
# Row 81
$ws.Range("H81").Value = 610.5
$ws.Range("I81").Value = 603.25
$ws.Range("J81").Value = 625
$ws.Range("K81").Value = 1206.5
$ws.Range("L81").Value = 1250
$ws.Range("M81").Value = -145.5
$ws.Range("N81").Value = -3372
# Row 84
$ws.Range("H84").Value = 610.5
$ws.Range("I84").Value = 603.25
$ws.Range("J84").Value = 625
$ws.Range("K84").Value = 6032.5
$ws.Range("L84").Value = 6250
$ws.Range("M84").Value = -728.5
$ws.Range("N84").Value = -16858
# Row 122
$ws.Range("H122").Value = 3537.8667
$ws.Range("I122").Value = 2958.625
$ws.Range("K122").Value = 8875.875
$ws.Range("M122").Value = -6425.875
# Row 132
$ws.Range("H132").Value = 3549537.2
$ws.Range("I132").Value = 5750275.5
$ws.Range("J132").Value = 3903.8333
$ws.Range("K132").Value = 17250826.5
$ws.Range("L132").Value = 11711.4999
$ws.Range("M132").Value = -17248296.5
$ws.Range("N132").Value = -16771.4999
# Row 133
$ws.Range("H133").Value = 52000
$ws.Range("J133").Value = 52000
$ws.Range("L133").Value = 52000
$ws.Range("N133").Value = -62120
# Row 136
$ws.Range("H136").Value = 10503963
$ws.Range("I136").Value = 11951992
$ws.Range("J136").Value = 5751.25
$ws.Range("K136").Value = 35855976
$ws.Range("L136").Value = 17253.75
$ws.Range("M136").Value = -35853426
$ws.Range("N136").Value = -22353.75
